$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6749
$ws.Range("H67").Value = 6749
$ws.Range("H74").Value = 5751.467
$ws.Range("I74").Value = 5642.6665
$ws.Range("J74").Value = 5914.6665
$ws.Range("K74").Value = 5642.6665
$ws.Range("L74").Value = 5914.6665
$ws.Range("M74").Value = -4706.6665
$ws.Range("N74").Value = -7786.6665
$ws.Range("H77").Value = 5751.467
$ws.Range("I77").Value = 5642.6665
$ws.Range("J77").Value = 5914.6665
$ws.Range("K77").Value = 28213.3325
$ws.Range("L77").Value = 29573.3325
$ws.Range("M77").Value = -23533.3325
$ws.Range("N77").Value = -38933.3325
$ws.Range("H80").Value = 71442870
$ws.Range("I80").Value = 142858450
$ws.Range("K80").Value = 428575350
$ws.Range("M80").Value = -428574352
$ws.Range("H83").Value = 71442870
$ws.Range("I83").Value = 142858450
$ws.Range("K83").Value = 1285726050
$ws.Range("M83").Value = -1285721058
$ws.Range("H88").Value = 2626.4
$ws.Range("I88").Value = 1709.3334
$ws.Range("J88").Value = 4002
$ws.Range("K88").Value = 1709.3334
$ws.Range("L88").Value = 4002
$ws.Range("M88").Value = -1303.3334
$ws.Range("N88").Value = -4814
$ws.Range("H91").Value = 2626.4
$ws.Range("I91").Value = 1709.3334
$ws.Range("J91").Value = 4002
$ws.Range("K91").Value = 1709.3334
$ws.Range("L91").Value = 4002
$ws.Range("M91").Value = -305.3334
$ws.Range("N91").Value = -6810
$ws.Range("H103").Value = 717.5172
$ws.Range("I103").Value = 556.5
$ws.Range("K103").Value = 1669.5
$ws.Range("M103").Value = -1083.5
$ws.Range("H111").Value = 3479.2
$ws.Range("I111").Value = 2754.6667
$ws.Range("K111").Value = 8264.000100000001
$ws.Range("M111").Value = -5197.000100000001
$ws.Range("H132").Value = 3019.238
$ws.Range("I132").Value = 2548.7334
$ws.Range("J132").Value = 4195.5
$ws.Range("K132").Value = 7646.2002
$ws.Range("L132").Value = 12586.5
$ws.Range("M132").Value = -5116.2002
$ws.Range("N132").Value = -17646.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2403.7273
$ws.Range("J2").Value = 2779.2856
$ws.Range("L2").Value = 2779.2856
$ws.Range("N2").Value = -3005.2856
$ws.Range("H32").Value = 19895442
$ws.Range("I32").Value = 20456100
$ws.Range("K32").Value = 20456100
$ws.Range("M32").Value = -20455813
$ws.Range("H45").Value = 3279.8635
$ws.Range("I45").Value = 2600.9
$ws.Range("K45").Value = 2600.9
$ws.Range("M45").Value = -2223.9
$ws.Range("H63").Value = 3812.3333
$ws.Range("I63").Value = 2417.2
$ws.Range("J63").Value = 4808.857
$ws.Range("K63").Value = 2417.2
$ws.Range("L63").Value = 4808.857
$ws.Range("M63").Value = -1731.2
$ws.Range("N63").Value = -6180.857
$ws.Range("H66").Value = 3812.3333
$ws.Range("I66").Value = 2417.2
$ws.Range("J66").Value = 4808.857
$ws.Range("K66").Value = 12086
$ws.Range("L66").Value = 24044.285
$ws.Range("M66").Value = -8654
$ws.Range("N66").Value = -30908.285
$ws.Range("H116").Value = 2403.7273
$ws.Range("J116").Value = 2779.2856
$ws.Range("L116").Value = 2779.2856
$ws.Range("N116").Value = -7367.2856
$ws.Range("H131").Value = 117923.5
$ws.Range("J131").Value = 117923.5
$ws.Range("L131").Value = 117923.5
$ws.Range("N131").Value = -128003.5
$ws.Range("H132").Value = 3019.3845
$ws.Range("I132").Value = 2619.35
$ws.Range("K132").Value = 7858.049999999999
$ws.Range("M132").Value = -5328.049999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2403.7273
$ws.Range("J3").Value = 2779.2856
$ws.Range("L3").Value = 2779.2856
$ws.Range("N3").Value = -3007.2856
$ws.Range("H99").Value = 2513.1428
$ws.Range("I99").Value = 2111
$ws.Range("J99").Value = 3800
$ws.Range("K99").Value = 2111
$ws.Range("L99").Value = 3800
$ws.Range("M99").Value = -613
$ws.Range("N99").Value = -6796

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H87").Value = 96883.5
$ws.Range("J87").Value = 96883.5
$ws.Range("L87").Value = 96883.5
$ws.Range("N87").Value = -99255.5
$ws.Range("H90").Value = 96883.5
$ws.Range("J90").Value = 96883.5
$ws.Range("L90").Value = 290650.5
$ws.Range("N90").Value = -302506.5
$ws.Range("H94").Value = 3515.125
$ws.Range("I94").Value = 5756
$ws.Range("J94").Value = 2768.1667
$ws.Range("K94").Value = 5756
$ws.Range("L94").Value = 2768.1667
$ws.Range("M94").Value = -5305
$ws.Range("N94").Value = -3670.1667
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 293
$ws.Range("I38").Value = 29.428572
$ws.Range("J38").Value = 908
$ws.Range("K38").Value = 88.28571599999999
$ws.Range("L38").Value = 2724
$ws.Range("M38").Value = 258.714284
$ws.Range("N38").Value = -3418
$ws.Range("H55").Value = 3852.7778
$ws.Range("J55").Value = 5004.1665
$ws.Range("L55").Value = 15012.4995
$ws.Range("N55").Value = -15366.4995
$ws.Range("H57").Value = 7120.6665
$ws.Range("J57").Value = 8504.799999999999
$ws.Range("L57").Value = 25514.4
$ws.Range("N57").Value = -26632.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1453.375
$ws.Range("J122").Value = 1999.6666
$ws.Range("L122").Value = 5998.9998
$ws.Range("N122").Value = -10898.9998
$ws.Range("H132").Value = 2928.7727
$ws.Range("I132").Value = 2924.1904
$ws.Range("K132").Value = 8772.5712
$ws.Range("M132").Value = -6242.5712

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7970
$ws.Range("I7").Value = 7912.4443
$ws.Range("J7").Value = 8073.6
$ws.Range("K7").Value = 7912.4443
$ws.Range("L7").Value = 8073.6
$ws.Range("M7").Value = -7800.4443
$ws.Range("N7").Value = -8297.6
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H100").Value = 1477.7142
$ws.Range("I100").Value = 1477.7142
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1477.7142
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -936.7141999999999
$ws.Range("N100").ClearContents()
$ws.Range("H116").Value = 158759.5
$ws.Range("J116").Value = 158759.5
$ws.Range("L116").Value = 158759.5
$ws.Range("N116").Value = -167937.5
$ws.Range("H122").Value = 13593.641
$ws.Range("I122").Value = 13701
$ws.Range("K122").Value = 41103
$ws.Range("M122").Value = -38653
$ws.Range("H126").Value = 7970
$ws.Range("I126").Value = 7912.4443
$ws.Range("J126").Value = 8073.6
$ws.Range("K126").Value = 23737.3329
$ws.Range("L126").Value = 24220.8
$ws.Range("M126").Value = -21267.3329
$ws.Range("N126").Value = -29160.8
$ws.Range("H132").Value = 5552.5415
$ws.Range("I132").Value = 5255.2856
$ws.Range("K132").Value = 15765.8568
$ws.Range("M132").Value = -13235.8568
$ws.Range("H136").Value = 10358.1
$ws.Range("I136").Value = 6697.625
$ws.Range("K136").Value = 20092.875
$ws.Range("M136").Value = -17542.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5366.6665
$ws.Range("J62").Value = 6425
$ws.Range("L62").Value = 6425
$ws.Range("N62").Value = -7673
$ws.Range("H65").Value = 5366.6665
$ws.Range("J65").Value = 6425
$ws.Range("L65").Value = 32125
$ws.Range("N65").Value = -38365
$ws.Range("H75").Value = 82473.8
$ws.Range("J75").Value = 87487.25
$ws.Range("L75").Value = 87487.25
$ws.Range("N75").Value = -89359.25
$ws.Range("H78").Value = 82473.8
$ws.Range("J78").Value = 87487.25
$ws.Range("L78").Value = 262461.75
$ws.Range("N78").Value = -271821.75
$ws.Range("H122").Value = 8097.7
$ws.Range("I122").Value = 8097.7
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 24293.1
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -21843.1
$ws.Range("N122").ClearContents()
$ws.Range("H123").Value = 113000
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()
$ws.Range("H132").Value = 2538.7693
$ws.Range("I132").Value = 2375.3333
$ws.Range("K132").Value = 7125.999899999999
$ws.Range("M132").Value = -4595.999899999999
